# Re-shuffle the per-observation data (species/coordinate columns) across
# rows 9-16 and 18 of the active sheet. Row 17 is left untouched.
#
# The columns that carry the "identity" of an observation record are:
#   A  (Id), B (Taxonsorteringsordning), D (Rödlistade), E (TaxonId),
#   F  (Artnamn), G (Vetenskapligt namn), H (Auktor), Q (Ost), R (Nord)
# All other columns (C, I, J, K, N, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AF, AG, AT, AW, AX, AY, ...) are identical for every one of these rows, so
# they do not need to move.
#
# The new row for each existing record (i.e. "row N after the edit shows
# the data that used to live on row M before the edit") is:
#   9 <- 12   10 <- 13   11 <- 10   12 <- 9   13 <- 16
#   14 <- 15  15 <- 14   16 <- 18   17 <- 17  18 <- 11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# destination row -> source row (values to copy FROM source row's current
# contents INTO the destination row)
$rowMap = [ordered]@{
    9  = 12
    10 = 13
    11 = 10
    12 = 9
    13 = 16
    14 = 15
    15 = 14
    16 = 18
    18 = 11
}

# 1) Snapshot the current ("before") values for every source row/column we
#    need, so that later writes don't clobber data we still have to read.
$snapshot = @{}
foreach ($srcRow in ($rowMap.Values | Sort-Object -Unique)) {
    foreach ($col in $cols) {
        $snapshot[$col + $srcRow] = $ws.Range($col + $srcRow).Value2
    }
}

# 2) Write the snapshotted values into their new destination rows.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Range($col + $destRow).Value2 = $snapshot[$col + $srcRow]
    }
}
